$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Prefix with an apostrophe so Excel stores numeric-looking strings
    # (e.g. "0.7121", "29.387.28") as literal text instead of coercing
    # them to a Number/Date, then reset the style so we don't leave a
    # stray "quote prefix" cell style behind.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "29.387.28"
$ws.Range("E2").Value = "  +0.54%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.874.57"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - XRP
Set-TextCell "D5" "0.7121"
$ws.Range("E5").Value = "  -0.28%  "

# Row 6 - BNB
Set-TextCell "D6" "242.11"
$ws.Range("E6").Value = "  +0.64%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - was Cardano, now Dogecoin
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell "D8" "0.07799"
$ws.Range("E8").Value = "  +1.07%  "

# Row 9 - was Dogecoin, now Cardano
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell "D9" "0.3111"
$ws.Range("E9").Value = "  +0.75%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  +0.82%  "

# Row 11 - TRON
Set-TextCell "D11" "0.08446"
$ws.Range("E11").Value = "  +1.62%  "

# Row 12 - WrappedEther
Set-TextCell "D12" "1.859.15"
$ws.Range("E12").Value = "  -0.57%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +0.50%  "

# Row 14 - Polygon
Set-TextCell "D14" "0.7125"
$ws.Range("E14").Value = "  -0.71%  "

# Row 15 - Litecoin
Set-TextCell "D15" "91.18"
$ws.Range("E15").Value = "  +0.33%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "29.388.34"
$ws.Range("E16").Value = "  +0.49%  "

# Row 17 - ShibaInu
Set-TextCell "D17" "0.000008233"
$ws.Range("E17").Value = "  +5.46%  "

# Row 18 - Uniswap
Set-TextCell "D18" "6.039"
$ws.Range("E18").Value = "  +0.75%  "

# Row 19 - BitcoinCash
Set-TextCell "D19" "241.16"
$ws.Range("E19").Value = "  -1.06%  "

# Row 20 - Avalanche
$ws.Range("E20").Value = "  +0.75%  "

# Row 21 - WrappedliquidstakedEther2.0
Set-TextCell "D21" "2.120.38"
$ws.Range("E21").Value = "  -0.87%  "

# Row 22 - Dai
Set-TextCell "D22" "0.9994"

# Row 23 - Chainlink
Set-TextCell "D23" "7.783"
$ws.Range("E23").Value = "  -2.16%  "

# Row 24 - BinanceUSD
Set-TextCell "D24" "1.001"
$ws.Range("E24").Value = "  +0.17%  "

# Row 25 - Stellar
Set-TextCell "D25" "0.1608"
$ws.Range("E25").Value = "  -0.24%  "

# Row 26 - Monero
Set-TextCell "D26" "163.55"
$ws.Range("E26").Value = "  +0.49%  "

# Row 27 - Cosmos
Set-TextCell "D27" "9.061"
$ws.Range("E27").Value = "  +1.63%  "

# Row 28 - EthereumClassic
Set-TextCell "D28" "18.48"
$ws.Range("E28").Value = "  -0.61%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.96%  "

# Row 30 - Filecoin
Set-TextCell "D30" "4.427"
$ws.Range("E30").Value = "  -0.25%  "

# Row 31 - Toncoin
Set-TextCell "D31" "1.289"
$ws.Range("E31").Value = "  -5.29%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextCell "D32" "4.307"
$ws.Range("E32").Value = "  +1.22%  "

# Row 33 - Hedera
Set-TextCell "D33" "0.05289"
$ws.Range("E33").Value = "  +2.01%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +0.43%  "

# Row 36 - ImmutableX
Set-TextCell "D36" "0.7442"
$ws.Range("E36").Value = "  -8.90%  "

# Row 37 - HuobiToken
Set-TextCell "D37" "2.699"
$ws.Range("E37").Value = "  +0.56%  "

# Row 38 - VeChain
Set-TextCell "D38" "0.01871"
$ws.Range("E38").Value = "  +0.64%  "

# Row 39 - Maker
Set-TextCell "D39" "1.214.62"
$ws.Range("E39").Value = "  +4.35%  "

# Row 40 - MXToken
Set-TextCell "D40" "2.724"
$ws.Range("E40").Value = "  +1.23%  "

# Row 41 - FraxShare
Set-TextCell "D41" "6.479"
$ws.Range("E41").Value = "  +4.29%  "

# Row 42 - TrustWalletToken
Set-TextCell "D42" "0.8879"
$ws.Range("E42").Value = "  -1.64%  "

# Row 43 - Aave
Set-TextCell "D43" "72.78"
$ws.Range("E43").Value = "  -0.12%  "

# Row 44 - Quant
Set-TextCell "D44" "108.98"
$ws.Range("E44").Value = "  +6.94%  "

# Row 45 - PaxDollar
$ws.Range("E45").Value = "  +0.11%  "

# Row 46 - RocketPoolETH
Set-TextCell "D46" "2.018.74"
$ws.Range("E46").Value = "  -1.26%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +1.84%  "

# Row 48 - Mantle
Set-TextCell "D48" "0.5210"
$ws.Range("E48").Value = "  +0.71%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  +2.52%  "

# Row 50 - EnergySwap
Set-TextCell "D50" "9.363"
$ws.Range("E50").Value = "  -0.22%  "

# Row 51 - TheSandbox
Set-TextCell "D51" "0.4324"
$ws.Range("E51").Value = "  +0.96%  "
